$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": append a new data row (row 50) after the last existing row.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

# Insert a fresh row 50 (this copies the formatting of row 49, which is NOT
# what the new row should look like for columns B:H - only column A's date
# style is reused). We then "launder" the formatting of B50:H50 back to the
# worksheet/column default by cutting it in from a pristine, never-touched
# range and pasting it on top - this avoids creating any new style entries.
$wsAll.Rows.Item(50).Insert()
$wsAll.Range("B200:H200").Cut($wsAll.Range("B50:H50"))

$wsAll.Range("A50").Value2 = 43978
$wsAll.Range("B50").Value2 = 285
$wsAll.Range("C50").Value2 = 282
$wsAll.Range("D50").Value2 = 18
$wsAll.Range("E50").Value2 = 15
$wsAll.Range("F50").Value2 = 3
$wsAll.Range("G50").Value2 = 12
$wsAll.Range("H50").Value2 = 252

# ---------------------------------------------------------------------------
# Sheet "kobe": insert a new data row (row 105) before the trailing note row,
# pushing the note row down to 106.
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Rows.Item(105).Insert()

$wsKobe.Range("A105").Value2 = 43978
$wsKobe.Range("B105").Value2 = 26
$wsKobe.Range("C105").Value2 = 3106
$wsKobe.Range("E105").Value2 = 285
$wsKobe.Range("F105").Value2 = 15
$wsKobe.Range("G105").Value2 = 13
$wsKobe.Range("H105").Value2 = 2
$wsKobe.Range("I105").Value2 = 12
$wsKobe.Range("J105").Value2 = 241

# ---------------------------------------------------------------------------
# Sheet "other": insert a new data row (row 80) before the trailing note row,
# pushing the note row down to 81.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

$wsOther.Rows.Item(80).Insert()

$wsOther.Range("A80").Value2 = 43978
$wsOther.Range("B80").Value2 = 0
$wsOther.Range("C80").Value2 = 14
$wsOther.Range("D80").Value2 = 3
$wsOther.Range("E80").Value2 = 2
$wsOther.Range("F80").Value2 = 1
$wsOther.Range("G80").Value2 = 0
$wsOther.Range("H80").Value2 = 11

# ---------------------------------------------------------------------------
# View state: the "kobe" sheet tab is now the active tab/sheet, and each
# sheet's saved selection moves to reflect where the editor ended up.
# ---------------------------------------------------------------------------
$wsAll.Activate()
$wsAll.Range("E52").Select()

$wsOther.Activate()
$wsOther.Range("D83").Select()

$wsKobe.Activate()
$wsKobe.Range("G108").Select()
